# Singapore Premier League workbook update
# The two fixtures played on 2023-06-06 (rows 4 & 5) had been recorded
# against the wrong teams/odds, and likewise the two fixtures played on
# 2023-09-19 (rows 54 & 55). This swaps each pair of rows back into their
# correct order (everything except the running index in column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $r1, $r2, $firstCol, $lastCol) {
    $rng1 = $ws.Range("$firstCol$r1`:$lastCol$r1")
    $rng2 = $ws.Range("$firstCol$r2`:$lastCol$r2")

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}

# Swap the full row content (columns B..AC) between rows 4 and 5,
# leaving column A (the sequential record index) untouched.
Swap-Rows $ws 4 5 "B" "AC"

# Swap the full row content (columns B..AC) between rows 54 and 55,
# leaving column A (the sequential record index) untouched.
Swap-Rows $ws 54 55 "B" "AC"
